$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 259, shifting existing rows 259:358 down to 261:360.
$ws.Rows("259:260").Insert()

# Populate the two newly inserted rows (259 and 260) with their data.

# Row 259
$ws.Range("A259").Value = 8
$ws.Range("B259").Value = "Terminal La Palmera de La Serena"
$ws.Range("C259").Value = "Coquimbo"
$ws.Range("D259").Value = 45006
$ws.Range("E259").Value = 4
$ws.Range("F259").Value = 100112031
$ws.Range("G259").Value = "Poroto verde"
$ws.Range("H259").Value = "Magnum"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 400
$ws.Range("K259").Value = 22000
$ws.Range("L259").Value = 23000
$ws.Range("M259").Value = 22500
$ws.Range("N259").Value = "$/caja 25 kilos"
$ws.Range("O259").Value = "Provincia de Limarí"
$ws.Range("P259").Value = 900
$ws.Range("Q259").Value = 25
$ws.Range("R259").Value = "Hortaliza"

# Row 260
$ws.Range("A260").Value = 8
$ws.Range("B260").Value = "Terminal La Palmera de La Serena"
$ws.Range("C260").Value = "Coquimbo"
$ws.Range("D260").Value = 45006
$ws.Range("E260").Value = 4
$ws.Range("F260").Value = 100112031
$ws.Range("G260").Value = "Poroto verde"
$ws.Range("H260").Value = "Magnum"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 400
$ws.Range("K260").Value = 23000
$ws.Range("L260").Value = 24000
$ws.Range("M260").Value = 23500
$ws.Range("N260").Value = "$/malla 25 kilos"
$ws.Range("O260").Value = "Provincia de Limarí"
$ws.Range("P260").Value = 940
$ws.Range("Q260").Value = 25
$ws.Range("R260").Value = "Hortaliza"

# Ensure date format on new D cells matches the rest of column D.
$ws.Range("D259:D260").NumberFormat = "YYYY-MM-DD HH:MM:SS"
